$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'32.927.51"
$ws.Range("E2").Value = "  +10.06%  "

# Row 3
$ws.Range("D3").Value = "'1.757.85"
$ws.Range("E3").Value = "  +5.99%  "

# Row 4
$ws.Range("D4").Value = "'0.996"
$ws.Range("E4").Value = "  -0.27%  "

# Row 5
$ws.Range("D5").Value = "'227.21"
$ws.Range("E5").Value = "  +4.74%  "

# Row 6
$ws.Range("D6").Value = "'0.543"
$ws.Range("E6").Value = "  +4.10%  "

# Row 7
$ws.Range("D7").Value = "'0.996"
$ws.Range("E7").Value = "  -0.22%  "

# Row 8
$ws.Range("D8").Value = "'31.94"
$ws.Range("E8").Value = "  +10.85%  "

# Row 9
$ws.Range("D9").Value = "'44.93"
$ws.Range("E9").Value = "  +2.53%  "

# Row 11
$ws.Range("D11").Value = "'0.0671"
$ws.Range("E11").Value = "  +9.15%  "

# Row 12
$ws.Range("D12").Value = "'0.0919"
$ws.Range("E12").Value = "  +1.82%  "

# Row 13
$ws.Range("D13").Value = "'2.009.83"
$ws.Range("E13").Value = "  +5.98%  "

# Row 14
$ws.Range("D14").Value = "'1.754.75"
$ws.Range("E14").Value = "  +5.81%  "

# Row 15
$ws.Range("D15").Value = "'0.632"
$ws.Range("E15").Value = "  +4.67%  "

# Row 16
$ws.Range("D16").Value = "'10.48"
$ws.Range("E16").Value = "  +5.01%  "

# Row 17
$ws.Range("D17").Value = "'4.29"
$ws.Range("E17").Value = "  +8.83%  "

# Row 18
$ws.Range("D18").Value = "'32.898.69"
$ws.Range("E18").Value = "  +9.85%  "

# Row 19
$ws.Range("D19").Value = "'68.91"
$ws.Range("E19").Value = "  +6.11%  "

# Row 20
$ws.Range("D20").Value = "'259.56"
$ws.Range("E20").Value = "  +7.14%  "

# Row 21
$ws.Range("D21").Value = "'0.0₃0742"
$ws.Range("E21").Value = "  +4.57%  "

# Row 22
$ws.Range("D22").Value = "'0.995"
$ws.Range("E22").Value = "  -0.35%  "

# Row 23
$ws.Range("D23").Value = "'10.54"
$ws.Range("E23").Value = "  +5.18%  "

# Row 24
$ws.Range("D24").Value = "'4.36"
$ws.Range("E24").Value = "  +4.41%  "

# Row 25
$ws.Range("D25").Value = "'2.18"
$ws.Range("E25").Value = "  -0.57%  "

# Row 26
$ws.Range("D26").Value = "'159.55"
$ws.Range("E26").Value = "  +0.67%  "

# Row 27
$ws.Range("D27").Value = "'16.56"
$ws.Range("E27").Value = "  +5.12%  "

# Row 28
$ws.Range("E28").Value = "  +4.28%  "

# Row 29
$ws.Range("D29").Value = "'6.98"
$ws.Range("E29").Value = "  +3.80%  "

# Row 30
$ws.Range("D30").Value = "'0.996"
$ws.Range("E30").Value = "  -0.23%  "

# Row 31
$ws.Range("E31").Value = "  +14.56%  "

# Row 32
$ws.Range("E32").Value = "  +3.49%  "

# Row 33
$ws.Range("E33").Value = "  +5.97%  "

# Row 34
$ws.Range("E34").Value = "  +7.91%  "

# Row 35
$ws.Range("D35").Value = "'1.550.94"
$ws.Range("E35").Value = "  +7.43%  "

# Row 36
$ws.Range("E36").Value = "  +4.83%  "

# Row 37
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "'1.05"
$ws.Range("E37").Value = "  +2.17%  "

# Row 38
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "'0.631"
$ws.Range("E38").Value = "  +9.92%  "

# Row 39
$ws.Range("D39").Value = "'84.82"
$ws.Range("E39").Value = "  +7.60%  "

# Row 40
$ws.Range("E40").Value = "  +6.22%  "

# Row 41
$ws.Range("E41").Value = "  +2.76%  "

# Row 42
$ws.Range("E42").Value = "  +0.24%  "

# Row 43
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "'0.875"
$ws.Range("E43").Value = "  +3.19%  "

# Row 44
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'2.09"
$ws.Range("E44").Value = "  +7.26%  "

# Row 45
$ws.Range("E45").Value = "  +2.50%  "

# Row 46
$ws.Range("D46").Value = "'54.76"
$ws.Range("E46").Value = "  +7.93%  "

# Row 47
$ws.Range("E47").Value = "  +4.54%  "

# Row 48
$ws.Range("D48").Value = "'1.908.87"
$ws.Range("E48").Value = "  +5.88%  "

# Row 49
$ws.Range("E49").Value = "  +6.06%  "

# Row 50
$ws.Range("D50").Value = "'0.997"
$ws.Range("E50").Value = "  -0.12%  "

# Row 51
$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").Value = "'95.85"
$ws.Range("E51").Value = "  +1.91%  "
